$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45178 -> 45179) for every data row (rows 2 through 292).
$ws.Range("C2:C292").Value = 45179
